$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 2370.8
$ws.Range("J2").Value = 2468.75
$ws.Range("L2").Value = 2468.75
$ws.Range("N2").Value = -2694.75
# Row 9
$ws.Range("H9").Value = 278.625
$ws.Range("I9").Value = 450
$ws.Range("J9").Value = 175.8
$ws.Range("K9").Value = 450
$ws.Range("L9").Value = 175.8
$ws.Range("M9").Value = -281
$ws.Range("N9").Value = -513.8
# Row 40
$ws.Range("H40").Value = 4600
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 5250
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 5250
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -5600
# Row 63
$ws.Range("H63").Value = 25271
$ws.Range("J63").Value = 25271
$ws.Range("L63").Value = 25271
$ws.Range("N63").Value = -26519
# Row 66
$ws.Range("H66").Value = 25271
$ws.Range("J66").Value = 25271
$ws.Range("L66").Value = 75813
$ws.Range("N66").Value = -82053
# Row 111
$ws.Range("H111").Value = 800
$ws.Range("J111").Value = 1000
$ws.Range("L111").Value = 3000
$ws.Range("N111").Value = -9134
# Row 132
$ws.Range("H132").Value = 1231.409
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
# Row 133
$ws.Range("H133").Value = 95512
$ws.Range("J133").Value = 95512
$ws.Range("L133").Value = 95512
$ws.Range("N133").Value = -105632

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 21783.846
$ws.Range("I32").Value = 4622.1753
$ws.Range("K32").Value = 4622.1753
$ws.Range("M32").Value = -4335.1753
# Row 45
$ws.Range("H45").Value = 5699.8
$ws.Range("I45").Value = 6987.8423
$ws.Range("K45").Value = 6987.8423
$ws.Range("M45").Value = -6610.8423
# Row 61
$ws.Range("H61").Value = 1718.1364
$ws.Range("I61").Value = 1689.95
$ws.Range("K61").Value = 1689.95
$ws.Range("M61").Value = -1477.95
# Row 97
$ws.Range("H97").Value = 630.58826
$ws.Range("I97").Value = 627.8125
$ws.Range("J97").Value = 675
$ws.Range("K97").Value = 627.8125
$ws.Range("L97").Value = 675
$ws.Range("M97").Value = -131.8125
$ws.Range("N97").Value = -1667
# Row 102
$ws.Range("H102").Value = 1861.8636
$ws.Range("I102").Value = 1122.75
$ws.Range("K102").Value = 1122.75
$ws.Range("M102").Value = 499.25
# Row 136
$ws.Range("H136").Value = 1718.1364
$ws.Range("I136").Value = 1689.95
$ws.Range("K136").Value = 5069.85
$ws.Range("M136").Value = -2519.85

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4320.108
$ws.Range("I20").Value = 3616.56
$ws.Range("J20").Value = 5785.8335
$ws.Range("K20").Value = 3616.56
$ws.Range("L20").Value = 5785.8335
$ws.Range("M20").Value = -3369.56
$ws.Range("N20").Value = -6279.8335
# Row 95
$ws.Range("H95").Value = 51560.75
$ws.Range("J95").Value = 51560.75
$ws.Range("L95").Value = 51560.75
$ws.Range("N95").Value = -57052.75
# Row 105
$ws.Range("H105").Value = 955.2174
$ws.Range("I105").Value = 951.6316
$ws.Range("K105").Value = 951.6316
$ws.Range("M105").Value = 795.3684

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 37
$ws.Range("H37").Value = 1500
$ws.Range("J37").Value = 1500
$ws.Range("L37").Value = 1500
$ws.Range("N37").Value = -1714
# Row 92
$ws.Range("H92").Value = 28100.5
$ws.Range("J92").Value = 28100.5
$ws.Range("L92").Value = 28100.5
$ws.Range("N92").Value = -33092.5
# Row 97
$ws.Range("H97").Value = 30197
$ws.Range("J97").Value = 30197
$ws.Range("L97").Value = 30197
$ws.Range("N97").Value = -32179
# Row 134
$ws.Range("H134").Value = 2868.4443
$ws.Range("I134").Value = 3014.0571
$ws.Range("K134").Value = 9042.1713
$ws.Range("M134").Value = -6507.1713

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 10539
$ws.Range("I56").Value = 10539
$ws.Range("K56").Value = 10539
$ws.Range("M56").Value = -10009
# Row 68
$ws.Range("H68").Value = 1757
$ws.Range("I68").Value = 1700
$ws.Range("K68").Value = 5100
$ws.Range("M68").Value = -4289
# Row 70
$ws.Range("H70").Value = 3554.3333
$ws.Range("I70").Value = 2613.3333
$ws.Range("J70").Value = 4495.3335
$ws.Range("K70").Value = 7839.999899999999
$ws.Range("L70").Value = 13486.0005
$ws.Range("M70").Value = -7524.999899999999
$ws.Range("N70").Value = -14116.0005
# Row 71
$ws.Range("H71").Value = 1757
$ws.Range("I71").Value = 1700
$ws.Range("K71").Value = 15300
$ws.Range("M71").Value = -11244
# Row 73
$ws.Range("H73").Value = 3554.3333
$ws.Range("I73").Value = 2613.3333
$ws.Range("J73").Value = 4495.3335
$ws.Range("K73").Value = 7839.999899999999
$ws.Range("L73").Value = 13486.0005
$ws.Range("M73").Value = -6747.999899999999
$ws.Range("N73").Value = -15670.0005

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 28282.666
$ws.Range("I20").Value = 2999.8572
$ws.Range("K20").Value = 2999.8572
$ws.Range("M20").Value = -2754.8572
# Row 92
$ws.Range("H92").Value = 1250
$ws.Range("J92").Value = 1250
$ws.Range("L92").Value = 1250
$ws.Range("N92").Value = -4994
# Row 113
$ws.Range("H113").Value = 4069
$ws.Range("I113").Value = 3748
$ws.Range("J113").Value = 5995
$ws.Range("K113").Value = 3748
$ws.Range("L113").Value = 5995
$ws.Range("M113").Value = -1578
$ws.Range("N113").Value = -10335
# Row 126
$ws.Range("H126").Value = 3353.4666
$ws.Range("I126").Value = 2758.1428
$ws.Range("J126").Value = 3874.375
$ws.Range("K126").Value = 8274.428400000001
$ws.Range("L126").Value = 11623.125
$ws.Range("M126").Value = -5804.428400000001
$ws.Range("N126").Value = -16563.125
# Row 132
$ws.Range("H132").Value = 5102.68
$ws.Range("I132").Value = 3977.2104
$ws.Range("J132").Value = 8666.666999999999
$ws.Range("K132").Value = 11931.6312
$ws.Range("L132").Value = 26000.001
$ws.Range("M132").Value = -9401.6312
$ws.Range("N132").Value = -31060.001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 56340.25
$ws.Range("I46").Value = 423123
$ws.Range("K46").Value = 423123
$ws.Range("M46").Value = -422935
# Row 61
$ws.Range("H61").Value = 98399.664
$ws.Range("I61").Value = 51413
$ws.Range("K61").Value = 51413
$ws.Range("M61").Value = -51211
# Row 96
$ws.Range("H96").Value = 28435
$ws.Range("J96").Value = 28435
$ws.Range("L96").Value = 28435
$ws.Range("N96").Value = -33927
# Row 113
$ws.Range("H113").Value = 98399.664
$ws.Range("I113").Value = 51413
$ws.Range("K113").Value = 51413
$ws.Range("M113").Value = -49243
# Row 123
$ws.Range("H123").Value = 57986.668
$ws.Range("J123").Value = 57986.668
$ws.Range("L123").Value = 57986.668
$ws.Range("N123").Value = -67786.66800000001
# Row 136
$ws.Range("H136").Value = 4241.3
$ws.Range("I136").Value = 3465.5908
$ws.Range("K136").Value = 10396.7724
$ws.Range("M136").Value = -7846.7724

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1297.1333
$ws.Range("I100").Value = 1050.6666
$ws.Range("J100").Value = 2283
$ws.Range("K100").Value = 2101.3332
$ws.Range("L100").Value = 4566
$ws.Range("M100").Value = -1560.3332
$ws.Range("N100").Value = -5648
# Row 122
$ws.Range("H122").Value = 1291.2142
$ws.Range("I122").Value = 1236.6923
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3710.0769
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1260.0769
$ws.Range("N122").Value = -10900
# Row 136
$ws.Range("H136").Value = 605
$ws.Range("I136").Value = 544.2308
$ws.Range("K136").Value = 1632.6924
$ws.Range("M136").Value = 917.3075999999999

Write-Host "Edit complete"
